$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.798.85'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.758.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.28%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5074'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.83%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.39'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.95%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2662'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06200'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.755.33'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.54%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06939'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.51%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6054'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.473'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.838.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006826'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.980.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.070'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.183'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.195'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '138.06'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.459'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.821'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.92%  '

$ws.Range("E29").Value = '  +5.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08222'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.688'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.405'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04373'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9991'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.653'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.84%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6054'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.731'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.47%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.940'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.89%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01544'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.53%  '

$ws.Range("E42").Value = '  -0.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3824'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.11%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7377'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.911'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05494'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1084'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.947'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '29.87'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.612'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.94%  '
